$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.224.37'
$ws.Range("E2").Value = '  +0.74%  '
$ws.Range("D3").Value = '1.605.67'
$ws.Range("E3").Value = '  +0.55%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.28'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.16%  '
$ws.Range("E6").Value = '  -0.10%  '
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.248'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("E9").Value = '  -0.39%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.22'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +1.55%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0813'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.30%  '
$ws.Range("D12").Value = '1.828.83'
$ws.Range("E12").Value = '  +0.53%  '
$ws.Range("D13").Value = '1.595.03'
$ws.Range("E13").Value = '  -0.11%  '
$ws.Range("E14").Value = '  +0.55%  '
$ws.Range("E15").Value = '  +0.41%  '
$ws.Range("D16").Value = '26.194.06'
$ws.Range("E16").Value = '  +0.70%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.86'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +2.52%  '
$ws.Range("E19").Value = '  -0.05%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '200.85'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -1.27%  '
$ws.Range("E21").Value = '  +0.93%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.28'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.21%  '
$ws.Range("E23").Value = '  +0.43%  '
$ws.Range("E24").Value = '  +1.98%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.19'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +2.11%  '
$ws.Range("E26").Value = '  -0.07%  '
$ws.Range("E27").Value = '  -2.58%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.20'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.06%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.56'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +1.78%  '
$ws.Range("E30").Value = '  +3.86%  '
$ws.Range("E31").Value = '  +0.36%  '
$ws.Range("E32").Value = '  +2.62%  '
$ws.Range("E33").Value = '  -0.76%  '
$ws.Range("E34").Value = '  +0.72%  '
$ws.Range("D36").Value = '1.166.04'
$ws.Range("E36").Value = '  +5.23%  '
$ws.Range("E37").Value = '  +3.81%  '
$ws.Range("E38").Value = '  -0.11%  '
$ws.Range("E39").Value = '  -0.26%  '
$ws.Range("E41").Value = '  +0.64%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.779'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.27%  '
$ws.Range("E43").Value = '  +3.97%  '
$ws.Range("D44").Value = '1.739.94'
$ws.Range("E44").Value = '  +0.45%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '91.56'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.96%  '
$ws.Range("E46").Value = '  +1.99%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '54.07'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +1.33%  '
$ws.Range("E48").Value = '  +0.25%  '
$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.407'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.31%  '
$ws.Range("B50").Value = 'USDD'
$ws.Range("C50").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.00'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.10%  '
$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").Value = '0.0₇0948'
$ws.Range("E51").Value = '  +2.51%  '
